$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.862.48"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "3.653.27"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'597.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.14%  "

$ws.Range("D6").Value = "'190.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.50%  "

$ws.Range("E7").Value = "  -1.46%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.699"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "'57.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.08%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.93%  "

$ws.Range("E12").Value = "  -6.45%  "

$ws.Range("D13").Value = "'10.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "

$ws.Range("D14").Value = "4.236.82"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("D15").Value = "3.655.47"
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "'18.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.18%  "

$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").Value = "67.631.93"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("E20").Value = "  -3.11%  "

$ws.Range("D21").Value = "'401.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("D22").Value = "'4.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").Value = "'87.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("D24").Value = "'11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").Value = "'2.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").Value = "'12.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.32%  "

$ws.Range("D27").Value = "'6.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "'3.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.90%  "

$ws.Range("D29").Value = "'9.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "

$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").Value = "'7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").Value = "'12.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'44.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.82%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'65.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("D36").Value = "'607.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "'0.393"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").Value = "0.0₃0769"
$ws.Range("E40").Value = "  -13.98%  "

$ws.Range("D41").Value = "'0.136"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  -3.42%  "

$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("E44").Value = "  -8.79%  "

$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("D46").Value = "2.775.88"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("D47").Value = "'3.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").Value = "'143.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.03%  "

$ws.Range("D49").Value = "'8.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.48%  "

$ws.Range("D50").Value = "'2.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.52%  "
